# origins_distinct_developers.xlsx edit script
# Rebuilds column A (Developer) with the revised, split-out developer list,
# fixes the two hyperlinked entries (Luden.io / pixel.lu) that move rows,
# clears/rewrites formatting on the affected cells, adds 9 new trailing rows,
# and sets column A's width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Full, alphabetically-sorted developer list (post-edit), one entry per
#    worksheet row starting at row 2 (row 1 is the "Developer" header).
# ---------------------------------------------------------------------------
$developers = @(
    '10tons',
    '2K Australia',
    '2K Boston (Irrational Games)',
    '2K Marin',
    '2pt Interactive',
    '3D Realms',
    '4A Games',
    '5th Cell',
    'Abstraction Games',
    'Acid Nerve',
    'ACQUIRE Corp',
    'Adglobe Live Wire',
    'Aggro Crab',
    'aheartfulofgames',
    'Airtight Games',
    'All Possible Futures',
    'Alphanu Game Studios',
    'Alvion',
    'Ant Workshop',
    'Applava',
    'Aqua Style',
    'Arc System Works',
    'Arkane Studios',
    'Arrowhead Game Studios',
    'Asobo Studio',
    'Asteroid Base',
    'Atari SA (Infogrames Lyon House)',
    'Atlus',
    'Atomicom',
    'Avalanche Software',
    'Awaceb',
    'Awfully Nice Studios',
    'B.B. Studio',
    'Bamtang Games',
    'Bandai Namco Studios',
    'Beam Team Games',
    'Beat Games',
    'Beautiful Glitch',
    'Beenox',
    'Behaviour Interactive (Artificial Mind and Movement)',
    'Ben Esposito',
    'Bend Studio (Eidetic)',
    'Benjamin Rivers Inc',
    'Bento Studio',
    'Bethesda Game Studios',
    'Big Ant Studios',
    'Big Huge Games',
    'Billy Goat Entertainment Ltd',
    'Bimboosoft',
    'BioWare',
    'Bit Planet Games',
    'Bithell Games',
    'Bizarre Creations',
    'Black Forest Games',
    'Black Matter Games',
    'Bloober Team',
    'Blue Isle Studios',
    'Blue Tongue Entertainment',
    'BlueTwelve Studio',
    'Bossa Studios',
    'Brainseed Factory',
    'Breaking Walls',
    'Brownies',
    'BUG-Studio',
    'Bugbear Entertainment',
    'Bungie',
    'Camouflaj',
    'Capcom',
    'Capcom Vancouver (Blue Castle Games)',
    'Capybara Games',
    'Cat Daddy Games',
    'Caustic Reality',
    'CCP Games',
    'CD Projekt Red',
    'Cellar Door Games',
    'Chequered Ink Ltd',
    'Christian Whitehead',
    'Chubby Pixel',
    'Clap Hanz',
    'Clever Beans',
    'Cloak and Dagger Games',
    'ClockStone Studios',
    'Clover Studio',
    'Codemasters',
    'Codemasters Birmingham',
    'Coffee Stain Studios',
    'Cold Beam Games',
    'Cold Iron Studios',
    'Cold Symmetry',
    'Coldwood Interactive',
    'Colossal Order',
    'ConcernedApe (Eric Barone)',
    'Counterplay Games',
    'Cowardly Creations',
    'Crave',
    'CrazyBunch',
    'CrazyLabs',
    'Criterion Games',
    'Crows Crows Crows',
    'Crystal Dynamics',
    'Crytek',
    'Curve Studios',
    'Cyanide',
    'Dambuster Studios (Deep Silver)',
    'Daylight Studios',
    'Dead Drop Studios LLC',
    'Deep Silver Volition',
    'Delphine Software International',
    'Dennaton Games',
    'Depth First Games',
    'Digital Eclipse',
    'Digital Extremes',
    'Digital Leisure Inc',
    'Dodge Roll',
    'Dontnod Entertainment',
    'Double Fine Productions',
    'DoubleMoose Games',
    'Dreadlocks Ltd',
    'Drinkbox Studios',
    'DryGin Studios',
    'Dynamighty',
    'EA Black Box',
    'EA DICE (Digital Illusions CE)',
    'EA Gothenburg (Ghost Games)',
    'EA Montreal',
    'EA Tiburon',
    'EA Vancouver',
    'Eat Sleep Play',
    'Eidos Montreal',
    'Ember Lab',
    'Empty Clip Studios',
    'Enhance Games',
    'Epic Games',
    'EQ-Games',
    'Eurocom Entertainment Software',
    'Evolution Studios',
    'Exact (EXcellent Application Create Team)',
    'Exordium Games',
    'Experiment 101',
    'Extra Mile Studios',
    'Facepalm Games',
    'FAKT Software',
    'Fallen Tree Games',
    'Fast Travel Games',
    'Fiddlesticks',
    'Firaxis Games',
    'FireForge Games',
    'Firesprite',
    'First Watch Games',
    'FitXR',
    'Flying Wild Hog',
    'Four Door Lemon',
    'Four Quarters',
    'Frame Interactive',
    'Fraoula',
    'Free Lives',
    'Frictional Games',
    'Frima Studio',
    'Frogwares',
    'FromSoftware',
    'Frontier Developments',
    'Frozenbyte',
    'Fun Bits Interactive',
    'Funbox Media',
    'Funselektor Labs Inc',
    'FuturLab',
    'Game Design Sweden AB',
    'Game Freak',
    'Game Swing',
    'Gameloft',
    'Gameloft Montreal',
    'Gaming Minds Studios',
    'Ganbarion',
    'Gearbox Quebec',
    'Gearbox Software',
    'Ghost Ship Games',
    'Ghost Town Games',
    'Giant Sparrow',
    'Giant Squid Studios',
    'GIANTS Software',
    'Gibier Games',
    'Glitchy Pixel',
    'Grab Games',
    'Graceful Decay',
    'Grasshopper Manufacture',
    'Grip Games',
    'Guerrilla Cambridge (SCE Studio Cambridge)',
    'Guerrilla Games',
    'Haemimont Games',
    'Hailstorm Games',
    'Halfbrick Studios',
    'Hangar 13 (2K Czech)',
    'HappyGiant',
    'Harmonix Music Systems',
    'Hazelight Studios',
    'HB Studios',
    'Heavy Iron Studios',
    'Hello Games',
    'Herobeat Studios',
    'High Impact Games',
    'High Moon Studios',
    'Hipster Whale',
    'Honey Parade Games',
    'Honeyslug',
    'Hopoo Games',
    'Hothead Games',
    'House House',
    'Housemarque',
    'Hucast Games',
    'HumaNature Studios',
    'HypeHype (Frogmind)',
    'I-Illusions',
    'id Software',
    'Image & Form',
    'Incognito Entertainment',
    'Infinity Ward',
    'Infogrames Sheffield House (Gremlin Interactive)',
    'Infuse Studio',
    'Insomniac Games',
    'Invisible Walls',
    'IO Interactive',
    'Ion Lands',
    'Ion Storm',
    'Ironwood Studios',
    'Isometricorp Games',
    'iSquared Games',
    'ISVR',
    'Ivory Tower',
    'IzHard',
    'Jammed Up Studios',
    'JanduSoft',
    'Johnson Voorsanger Productions',
    'Juggler Games',
    'Just Add Water',
    'Jutsu Games',
    'Klabater',
    'Klei Entertainment',
    'Konami',
    'Krisalis Software',
    'Kronos Digital Entertainment',
    'Kung Fu Factory',
    'Kylotonn',
    'Landfall Games',
    'Laser Dog',
    'Laughing Jackal Ltd',
    'Le Cartel Studio',
    'League of Geeks',
    'Level-5',
    'Lichthund',
    'Lightwood Games',
    'Liquid Pug',
    'Little Chicken',
    'Lizardcube',
    'Llamasoft',
    'London Studio',
    'Louis Rigaud',
    'LucasArts',
    'Lucid Games',
    'Luden.io',
    'Ludosity',
    'MachineGames',
    'Maddy Makes Games',
    'maJAJa',
    'Marvelous Interactive',
    'Masaya',
    'Mass Creation',
    'MassHive Media',
    'Massive Monster',
    'Maxis',
    'Media Molecule',
    'Mediatonic',
    'Mega Crit',
    'metricminds',
    'Midway',
    'Mighty Rocket Studio',
    'miHoYo',
    'Milestone',
    'MixedBag',
    'Modern Dream',
    'Mojang Studios',
    'Monolith Productions',
    'Monstars',
    'Moppin (Ojiro Fumoto)',
    'Mossmouth LLC',
    'Mothership Entertainment',
    'Motive Studio',
    'Mouldy Toof Studios',
    'Nadeo',
    'NanaOn-Sha',
    'Naughty Dog',
    'nDreams',
    'Neat Corporation',
    'Necrophone Games',
    'NetherRealm Studios',
    'Nicalis',
    'Night School Studio',
    'Nihilistic Software',
    'Nippon Ichi Software',
    'No Brakes Games',
    'NomNomNami',
    'Noowanda',
    'Norsfell Games',
    'Nurijoy',
    'Obsidian Entertainment',
    'Odd Bug Studios',
    'Oddboy',
    'Oddworld Inhabitants',
    'Omega Force',
    'One Loop Games',
    'Other Ocean Emeryville',
    'Other Ocean Interactive',
    'Overkill Software',
    'Owlchemy Labs',
    'Pandemic Studios',
    'Paon',
    'Parabole',
    'Paradigm Entertainment',
    'Passtech Games',
    'Pencil Test Studios',
    'Pendulo Studios',
    'Perfectly Paranormal',
    'Pinokl Games',
    'Pixel Crow',
    'Pixel Maniacs',
    'pixel.lu',
    'Pixelopus',
    'PlatinumGames',
    'Play-Em',
    'Playbox (Digital Reality)',
    'Playrise Digital Ltd',
    'Playtonic Games',
    'Polyarc',
    'Polyphony Digital',
    'Polytron Corporation',
    'PopCap Games',
    'Pretty Fly Studios',
    'Prospect Games',
    'Pseudo Interactive',
    'Psyonix',
    'Puppy Games',
    'Purple Lamp Studios',
    'Q-Games',
    'Quantic Dream',
    'Queasy Games',
    'Radical Entertainment',
    'RageSquid',
    'Rampage Game Studios',
    'Ratloop',
    'Ratloop Asia',
    'Re-Logic',
    'Ready At Dawn Studios',
    'Realmforge Studios',
    'Rebellion Developments (Core Design)',
    'Red Barrels',
    'Red Phantom Games',
    'Red Storm Entertainment',
    'Reddoll Srl',
    'RedLynx',
    'Relentless Software',
    'Remedy Entertainment',
    'Repixel8',
    'Respawn Entertainment',
    'Robomodo',
    'Robot Gentleman',
    'Rockstar Leeds',
    'Rockstar North',
    'Rockstar San Diego',
    'Rockstar Toronto',
    'Rockstar Vancouver',
    'Rocksteady Studios',
    'Roll7',
    'Romero Games',
    'Rovio Entertainment (Housemarque)',
    'Running With Scissors',
    'Ryu Ga Gotoku Studio',
    'Saber Interactive',
    'Samurai Punk',
    'San Diego Studio',
    'Santa Monica Studio',
    'Sanzaru Games',
    'Schell Games',
    'Secret Sorcery',
    'SEGA',
    'Sergey Noskov',
    'Shared Memory',
    'Shedworks',
    'Shin`en Multimedia',
    'SingleTrac',
    'Skookum Arts',
    'Skydance Interactive',
    'Slick Entertainment',
    'Slightly Mad Studios',
    'Slipgate Ironworks',
    'Sloclap',
    'SMG Studio',
    'Snoozy Kazoo',
    'SoMa Play',
    'Something We Made',
    'Sometimes You',
    'Sony Interactive Entertainment',
    'Sony Interactive Studios America (989 Studios)',
    'SouthPAW Games',
    'souvenir circ',
    'Spicy Horse',
    'Spiders',
    'Spike Chunsoft',
    'Spiral House',
    'Splash Damage',
    'Spooky Doorway',
    'Sports Interactive',
    'Squad',
    'Squanch Games',
    'Square Enix (SquareSoft)',
    'Starbreeze Studios',
    'Steel Crate Games',
    'Striking Distance Studios',
    'Studio Liverpool (Psygnosis)',
    'Studio Wildcard',
    'Sucker Punch Productions',
    'Sumo Digital',
    'Superbot Entertainment',
    'Superflat Games',
    'Superhot Team',
    'Supermassive Games',
    'Supersonic Software',
    'Survios',
    'System 3 Software',
    'Systemic Reaction',
    'Tabot Inc',
    'Tamsoft',
    'Tango Gameworks',
    'Tarsier Studios',
    'Team Arcana (Examu)',
    'Team Asobi (Japan Studio)',
    'Team Bondi',
    'Team Cherry',
    'Team Meat',
    'Team Ninja',
    'Team17',
    'Telltale Games',
    'Tequila Works',
    'Terminal Reality',
    'Terrifying Jellyfish',
    'Teyon',
    'Thatgamecompany',
    'The Behemoth',
    'The Binary Mill',
    'The Chinese Room',
    'The Domaginarium',
    'The Fox Software',
    'The Game Atelier',
    'The Molasses Flood',
    'The Munky',
    'The Wild Gentlemen',
    'Thekla Inc',
    'Thinice Games',
    'Thomas Happ Games',
    'Three Fields Entertainment',
    'Tiger & Squid',
    'Tikipod',
    'Titan Studios',
    'Toby Fox',
    'Toge Productions',
    'Toys For Bob',
    'Treyarch',
    'Tribute Games (Jonathan Lavigne)',
    'Trinity Team',
    'Triple Eh? Ltd',
    'Tripwire Interactive',
    'TT Games (Travellers Tales)',
    'Turbo Button',
    'Twistplay',
    'Two Point Studios',
    'Ubisoft Chengdu',
    'Ubisoft Leamington (FreeStyleGames and Exient Entertainment)',
    'Ubisoft Montpellier',
    'Ubisoft Montreal',
    'Ubisoft Reflections (Reflections Interactive)',
    'Ubisoft San Francisco',
    'Ubisoft Shanghai',
    'Ubisoft Toronto',
    'Unbroken Studios',
    'Uncommon Chocolate',
    'Unfinished Pixel',
    'United Front Games',
    'Unknown Worlds Entertainment',
    'Upfall Studios',
    'Vblank Entertainment',
    'Vector Unit',
    'Velan Studios',
    'Vertigo Games',
    'Vicarious Visions',
    'Viewpoint Games',
    'Vile Monarch',
    'Visceral Games (EA Redwood Shores)',
    'Visual Concepts',
    'VitruviusVR',
    'VooFoo Studios',
    'VR Factory',
    'VRMonkey',
    'Wales Interactive',
    'WayForward Technologies Inc',
    'WayTooManyGames',
    'Weappy Studio',
    'WeirdBeard',
    'Whirlybird Games',
    'Whoopee Camp',
    'Wildbit Studios',
    'Wish Studios',
    'Wishes Unlimited (Greg Lobanov)',
    'Witch Beam',
    'Wolf & Wood Interactive Ltd',
    'Wolf Brew Games',
    'WolfEye Studios',
    'Workyrie Game Studio',
    'Xaloc Studios',
    'Yacht Club Games',
    'Yager Development',
    'Yeah Us!',
    'Young Horses',
    'Yukes'
)

$hyperlinkNames = @{
    'Luden.io' = 'http://luden.io'
    'pixel.lu' = 'http://pixel.lu'
}

# Remove every existing hyperlink up front -- the engine's Hyperlinks.Delete()
# is worksheet-scoped, so this always clears both of the old ones regardless
# of which range it's called on.
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Range("A1").Hyperlinks.Delete()
}

for ($i = 0; $i -lt $developers.Length; $i++) {
    $row = $i + 2
    $name = $developers[$i]
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $name

    if ($hyperlinkNames.ContainsKey($name)) {
        # Hyperlinked developer name: single underline, blue font.
        $cell.Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle
        $cell.Font.Color = 16711680
        $cell.HorizontalAlignment = [Microsoft.Office.Interop.Excel.Constants]::xlCenter
        $ws.Hyperlinks.Add($cell, $hyperlinkNames[$name]) | Out-Null
    } else {
        # Regular developer name: plain font, centered, no underline.
        $cell.Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleNone
        $cell.Font.Color = 0
        $cell.HorizontalAlignment = [Microsoft.Office.Interop.Excel.Constants]::xlCenter
    }
}

# ---------------------------------------------------------------------------
# 2) Column A width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 52
